{"js": "// Apply the day-of-week/date header update plus all 25 division-problem\n// text replacements described by the diff. Every \"old\" value occurs exactly\n// once in the document and every \"new\" value is unique, so a straight\n// search-and-replace per pair is safe regardless of execution order.\nconst replacements = [\n  [\"2025-02-22 Saturday\", \"2025-02-23 Sunday\"],\n  [\"220\u00f74=\", \"945\u00f76=\"],\n  [\"771\u00f73=\", \"531\u00f75=\"],\n  [\"332\u00f74=\", \"725\u00f75=\"],\n  [\"331\u00f75=\", \"104\u00f72=\"],\n  [\"391\u00f76=\", \"536\u00f75=\"],\n  [\"236\u00f74=\", \"160\u00f73=\"],\n  [\"247\u00f73=\", \"904\u00f74=\"],\n  [\"533\u00f76=\", \"909\u00f74=\"],\n  [\"527\u00f78=\", \"575\u00f78=\"],\n  [\"525\u00f79=\", \"768\u00f78=\"],\n  [\"194\u00f76=\", \"705\u00f79=\"],\n  [\"992\u00f79=\", \"624\u00f74=\"],\n  [\"737\u00f76=\", \"570\u00f79=\"],\n  [\"978\u00f73=\", \"620\u00f76=\"],\n  [\"907\u00f76=\", \"878\u00f79=\"],\n  [\"758\u00f74=\", \"433\u00f78=\"],\n  [\"262\u00f78=\", \"701\u00f76=\"],\n  [\"442\u00f77=\", \"941\u00f73=\"],\n  [\"144\u00f75=\", \"311\u00f73=\"],\n  [\"513\u00f74=\", \"554\u00f79=\"],\n  [\"543\u00f75=\", \"481\u00f77=\"],\n  [\"229\u00f76=\", \"948\u00f79=\"],\n  [\"359\u00f73=\", \"418\u00f74=\"],\n  [\"885\u00f76=\", \"863\u00f79=\"],\n  [\"578\u00f77=\", \"282\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the day-of-week/date header update plus all 25 division-problem\n# text replacements described by the diff. Every \"Old\" value occurs exactly\n# once in the document and every \"New\" value is unique, so a straight\n# find/replace per pair is safe regardless of execution order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-02-22 Saturday\"; New = \"2025-02-23 Sunday\" },\n    @{ Old = \"220\u00f74=\"; New = \"945\u00f76=\" },\n    @{ Old = \"771\u00f73=\"; New = \"531\u00f75=\" },\n    @{ Old = \"332\u00f74=\"; New = \"725\u00f75=\" },\n    @{ Old = \"331\u00f75=\"; New = \"104\u00f72=\" },\n    @{ Old = \"391\u00f76=\"; New = \"536\u00f75=\" },\n    @{ Old = \"236\u00f74=\"; New = \"160\u00f73=\" },\n    @{ Old = \"247\u00f73=\"; New = \"904\u00f74=\" },\n    @{ Old = \"533\u00f76=\"; New = \"909\u00f74=\" },\n    @{ Old = \"527\u00f78=\"; New = \"575\u00f78=\" },\n    @{ Old = \"525\u00f79=\"; New = \"768\u00f78=\" },\n    @{ Old = \"194\u00f76=\"; New = \"705\u00f79=\" },\n    @{ Old = \"992\u00f79=\"; New = \"624\u00f74=\" },\n    @{ Old = \"737\u00f76=\"; New = \"570\u00f79=\" },\n    @{ Old = \"978\u00f73=\"; New = \"620\u00f76=\" },\n    @{ Old = \"907\u00f76=\"; New = \"878\u00f79=\" },\n    @{ Old = \"758\u00f74=\"; New = \"433\u00f78=\" },\n    @{ Old = \"262\u00f78=\"; New = \"701\u00f76=\" },\n    @{ Old = \"442\u00f77=\"; New = \"941\u00f73=\" },\n    @{ Old = \"144\u00f75=\"; New = \"311\u00f73=\" },\n    @{ Old = \"513\u00f74=\"; New = \"554\u00f79=\" },\n    @{ Old = \"543\u00f75=\"; New = \"481\u00f77=\" },\n    @{ Old = \"229\u00f76=\"; New = \"948\u00f79=\" },\n    @{ Old = \"359\u00f73=\"; New = \"418\u00f74=\" },\n    @{ Old = \"885\u00f76=\"; New = \"863\u00f79=\" },\n    @{ Old = \"578\u00f77=\"; New = \"282\u00f73=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
